$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A15").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A16").Value = '{''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A17").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A18").Value = '{''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A19").Value = '{''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A21").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A22").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A23").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A24").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A30").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A31").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A32").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A33").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A35").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A36").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A37").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A38").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A39").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A40").Value = '{''FormalEducation'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A41").Value = '{''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A42").Value = '{''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A43").Value = '{''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A44").Value = '{''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A46").Value = '{''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A47").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A48").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A49").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A50").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A51").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A53").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A54").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A55").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A57").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A58").Value = '{''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A60").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1''}'
$ws.Range("A62").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1''}'
$ws.Range("A64").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A65").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A66").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A67").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A69").Value = '{''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A70").Value = '{''Hobby'': ''1'', ''DevType'': ''2''}'
$ws.Range("A72").Value = '{''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A73").Value = '{''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A74").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A75").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A77").Value = '{''HDI'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A78").Value = '{''GINI'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A79").Value = '{''GINI'': ''2'', ''Student'': ''1''}'
$ws.Range("A81").Value = '{''Age'': ''3'', ''SexualOrientation'': ''1''}'
$ws.Range("A82").Value = '{''Age'': ''3'', ''Student'': ''1''}'
$ws.Range("A83").Value = '{''Age'': ''3'', ''Hobby'': ''1''}'
